# Update the ClipOffset values (column E) on Sheet1 with new offset times
# pulled from PremPro. Column F ("HazardWindow") recalculates automatically
# since it holds the formula =E-D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$offsets = @{
    2  = 25.3
    3  = 37.616666666666667
    4  = 32.4
    5  = 19.166666666666668
    6  = 16.283333333333335
    7  = 33.966666666666669
    8  = 40.466666666666669
    9  = 46.983333333333334
    10 = 20.133333333333333
    11 = 36.06666666666667
    12 = 26.383333333333333
    13 = 29.65
    14 = 25.116666666666667
    15 = 12.9
    16 = 20.566666666666666
    17 = 15.633333333333333
    18 = 16.383333333333333
    19 = 33.583333333333336
    20 = 27.65
    21 = 13.8
    22 = 28.833333333333332
    23 = 25.633333333333333
    24 = 18.05
    25 = 39.583333333333336
}

foreach ($row in ($offsets.Keys | Sort-Object)) {
    $ws.Cells.Item($row, 5).Value = $offsets[$row]
}

$ws.Activate()
$ws.Range("F4").Select()
